# add DMSO-treated mice data from Kostourou et al., 2012
# New row is inserted into the "Vessel density (tumor)" sheet (4th sheet),
# right after the existing "Kostourou et al., 2013 (C57BL6/129 & CMT19T cell)"
# row and before the "Jones et al., 2013 (C57BL6 & LLC cell)" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Push the existing rows 9-11 down to make room for the new data point.
$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = "Kostourou et al., 2013 (C57BL6/129 & B16F0 cell & DMSO)"
$ws.Range("B9").Value = 64.849999999999994
$ws.Range("C9").Value = 7

# Grow Table4 so the new row is included in the table range.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C12"))

# Reflect the selection/active-sheet state captured in the saved workbook.
$null = $ws.Activate()
$null = $ws.Range("A10").Select()
